# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (placed right after "总计") with the
# quarter's fund-holdings detail, and adds a matching summary row at the
# top of "总计"'s data table (pushing the existing quarters down by one
# row).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New "2022-Q3" worksheet, right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row index column (numeric)
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(5,1).Value = 3

# Numeric-looking columns (fund code / size / position figures) must stay
# TEXT, matching the source data -- force with a temporary "@" format so
# Excel doesn't silently coerce them to numbers, then drop the format
# again so the cell keeps the workbook's default (unstyled) look.
$textCells = @(
    $q3.Cells.Item(2,2), $q3.Cells.Item(2,4), $q3.Cells.Item(2,5), $q3.Cells.Item(2,6), $q3.Cells.Item(2,7),
    $q3.Cells.Item(3,2), $q3.Cells.Item(3,4), $q3.Cells.Item(3,5), $q3.Cells.Item(3,6), $q3.Cells.Item(3,7),
    $q3.Cells.Item(4,2), $q3.Cells.Item(4,4), $q3.Cells.Item(4,5), $q3.Cells.Item(4,6), $q3.Cells.Item(4,7),
    $q3.Cells.Item(5,2), $q3.Cells.Item(5,4), $q3.Cells.Item(5,5), $q3.Cells.Item(5,6), $q3.Cells.Item(5,7)
)
foreach ($c in $textCells) { $c.NumberFormat = "@" }

$q3.Cells.Item(2,2).Value = "206013"
$q3.Cells.Item(2,3).Value = "鹏华宏观灵活配置混合"
$q3.Cells.Item(2,4).Value = "0.95"
$q3.Cells.Item(2,5).Value = "72.10"
$q3.Cells.Item(2,6).Value = "3.65"
$q3.Cells.Item(2,7).Value = "0.0347"
$q3.Cells.Item(2,8).Value = 9

$q3.Cells.Item(3,2).Value = "970042"
$q3.Cells.Item(3,3).Value = "国海量化优选一年持有股票C"
$q3.Cells.Item(3,4).Value = "7.16"
$q3.Cells.Item(3,5).Value = "87.31"
$q3.Cells.Item(3,6).Value = "0.34"
$q3.Cells.Item(3,7).Value = "0.0243"
$q3.Cells.Item(3,8).Value = 8

$q3.Cells.Item(4,2).Value = "562530"
$q3.Cells.Item(4,3).Value = "华夏中证智选1000价值稳健策略ETF"
$q3.Cells.Item(4,4).Value = "0.54"
$q3.Cells.Item(4,5).Value = "94.32"
$q3.Cells.Item(4,6).Value = "0.89"
$q3.Cells.Item(4,7).Value = "0.0048"
$q3.Cells.Item(4,8).Value = 10

$q3.Cells.Item(5,2).Value = "970041"
$q3.Cells.Item(5,3).Value = "国海量化优选一年持有股票A"
$q3.Cells.Item(5,4).Value = "0.63"
$q3.Cells.Item(5,5).Value = "87.31"
$q3.Cells.Item(5,6).Value = "0.34"
$q3.Cells.Item(5,7).Value = "0.0021"
$q3.Cells.Item(5,8).Value = 8

foreach ($c in $textCells) {
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Match the header/index-column formatting used on the other sheets
# (bold, centered, thin border == totalSheet's "B1" / "A2" style).
$totalSheet.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$totalSheet.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. Insert a "2022-Q3" summary row at the top of "总计"'s data
#    (row 2), pushing the other quarters down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q3"
$totalSheet.Cells.Item(2,3).Value = 4
$totalSheet.Cells.Item(2,4).Value = 0.07000000000000001

# Re-apply the index-column style to the new row (row-insert pulls in
# a slightly different blended format otherwise).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
$totalSheet.Cells.Item(2,1).Value = 0

# Renumber the index column (A) for the rows that shifted down.
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3

# Restore "总计" as the active sheet/selection (unchanged by this edit).
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
